$d = $word.ActiveDocument

function Merge-RunsAt($anchorText) {
    # Locate a unique, stable anchor string and force Word to re-coalesce
    # adjacent runs that share identical formatting within that paragraph.
    # We do this by toggling the very first character of the anchor to a
    # placeholder value and then restoring it — the structural edit this
    # triggers causes the engine to merge runs with identical rPr that sit
    # next to each other, which is exactly what the target diff needs
    # (several runs holding parts of one sentence get merged into a single
    # run once the redundant run-split is removed).
    $full = $d.Content
    $found = $full.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Anchor not found: $anchorText"
    }
    $p0 = $full.Start
    $orig = $d.Range($p0, $p0 + 1).Text

    $tmp = $d.Range($p0, $p0 + 1)
    $tmp.Text = "@"

    $restore = $d.Range($p0, $p0 + 1)
    $restore.Text = $orig
}

# Hunk 1: "(button) " + "will navigate to sign in page." -> merged run
Merge-RunsAt "Sign in"

# Hunk 2: "(link) " + "will navigate to sign up page." -> merged run
Merge-RunsAt "Sign up "

# Hunk 3: "(link) " + "will navigate to about us section in Home page (current page)." -> merged run
Merge-RunsAt "About us "

# Hunk 4: "(link) " + "will navigate to job applying section in Home page (current page)." -> merged run
Merge-RunsAt "Join our team "

# Hunk 5 + Hunk 6 (same paragraph):
#   "- " + "Small section of pictures ... wolt website - " -> merged run
#   "." + " " (after the wolt.com hyperlink) -> merged run
Merge-RunsAt "About us section"

# Hunk 7: "- (button) " + "will navigate " + "work with us page" + "." -> merged run
Merge-RunsAt "Check out our avaliable jobs "

Write-Output "done"
